$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (new most-recent reporting period).
# This shifts existing D..K one column right to E..L.
$ws.Range("D:D").EntireColumn.Insert()

# Copy cell formatting (number formats/styles) from the shifted column E
# into the new column D so each row keeps the same look (date header rows,
# number rows, etc.). Done per contiguous block so we don't touch the
# label-only separator rows (36/37/78/79) that have no data columns.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 112700
$ws.Range("D9").Value = 71700
$ws.Range("D10").Value = 41000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 800
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 81300
$ws.Range("D18").Value = 31400
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 31600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 31400
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 31300
$ws.Range("D27").Value = 28000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 28000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 28000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 55200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 14600
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 7700
$ws.Range("D46").Value = 77400
$ws.Range("D47").Value = 1507800
$ws.Range("D48").Value = 17600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1602900
$ws.Range("D57").Value = 6800
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 6800
$ws.Range("D61").Value = 727700
$ws.Range("D62").Value = 534100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1302000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 41100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 300800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 28000
$ws.Range("D83").Value = 200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 200
$ws.Range("D91").Value = "NA"
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -190400
$ws.Range("D96").Value = -22900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 164600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -25600
$ws.Range("E94").Value = -345900
$ws.Range("F94").Value = -299200
$ws.Range("E102").Value = 43800
